$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (old D:K -> new E:L)
$ws.Columns.Item(4).Insert()

# Copy number formats/styles from the (now shifted) neighbor column E into the
# new column D so every cell in D7:D102 gets the same style its row already
# uses (style 2 for the three "Period Ending" date rows, style 3 elsewhere).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

# Populate the new column D with the newly reported period's figures.
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(8, 4).Value = 14435000
$ws.Cells.Item(9, 4).Value = 11702000
$ws.Cells.Item(10, 4).Value = 2733000
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(14, 4).Value = 157000
$ws.Cells.Item(15, 4).Value = 124000
$ws.Cells.Item(17, 4).Value = 12976000
$ws.Cells.Item(18, 4).Value = 1459000
$ws.Cells.Item(20, 4).Value = 16000
$ws.Cells.Item(21, 4).Value = 2151000
$ws.Cells.Item(22, 4).Value = 141000
$ws.Cells.Item(23, 4).Value = 1334000
$ws.Cells.Item(24, 4).Value = 220000
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(26, 4).Value = 1114000
$ws.Cells.Item(27, 4).Value = 1097000
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(29, 4).Value = -30000
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(32, 4).Value = -16000
$ws.Cells.Item(33, 4).Value = 1067000
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(35, 4).Value = 1067000
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(41, 4).Value = 567000
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(43, 4).Value = 2795000
$ws.Cells.Item(44, 4).Value = 1277000
$ws.Cells.Item(45, 4).Value = 138000
$ws.Cells.Item(46, 4).Value = 4777000
$ws.Cells.Item(47, 4).Value = 215000
$ws.Cells.Item(48, 4).Value = 3179000
$ws.Cells.Item(49, 4).Value = 3904000
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(52, 4).Value = 405000
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(54, 4).Value = 12480000
$ws.Cells.Item(57, 4).Value = 2334000
$ws.Cells.Item(58, 4).Value = 306000
$ws.Cells.Item(59, 4).Value = 1054000
$ws.Cells.Item(60, 4).Value = 3694000
$ws.Cells.Item(61, 4).Value = 4038000
$ws.Cells.Item(62, 4).Value = 1078000
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(66, 4).Value = 9021000
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(72, 4).Value = 2511000
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(76, 4).Value = 3459000
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(81, 4).Value = 1067000
$ws.Cells.Item(83, 4).Value = 676000
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(89, 4).Value = 1628000
$ws.Cells.Item(91, 4).Value = -846000
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(94, 4).Value = -2048000
$ws.Cells.Item(96, 4).Value = -233000
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(100, 4).Value = -555000
$ws.Cells.Item(101, 4).Value = -54000
$ws.Cells.Item(102, 4).Value = -1029000
